# Fruta / hortaliza, semanal
# Insert a new week's pair of rows (Cebolla, Terminal La Palmera de La Serena)
# right before row 447, pushing the existing rows 447:472 down to 449:474.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 447; this shifts rows 447:472 -> 449:474
$ws.Rows("447:448").Insert()

# New row 447 data (1a (cosecha))
$ws.Cells.Item(447, 1).Value = 8
$ws.Cells.Item(447, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(447, 3).Value = "Coquimbo"
$ws.Cells.Item(447, 4).Value = 44516
$ws.Cells.Item(447, 4).NumberFormat = $ws.Cells.Item(449, 4).NumberFormat
$ws.Cells.Item(447, 5).Value = 4
$ws.Cells.Item(447, 6).Value = 100112004
$ws.Cells.Item(447, 7).Value = "Cebolla"
$ws.Cells.Item(447, 8).Value = "Sin especificar"
$ws.Cells.Item(447, 9).Value = "1a (cosecha)"
$ws.Cells.Item(447, 10).Value = 3200
$ws.Cells.Item(447, 11).Value = 4800
$ws.Cells.Item(447, 12).Value = 5000
$ws.Cells.Item(447, 13).Value = 4900
$ws.Cells.Item(447, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(447, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(447, 16).Value = 272
$ws.Cells.Item(447, 17).Value = 18
$ws.Cells.Item(447, 18).Value = "Hortaliza"

# New row 448 data (2a (cosecha))
$ws.Cells.Item(448, 1).Value = 8
$ws.Cells.Item(448, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(448, 3).Value = "Coquimbo"
$ws.Cells.Item(448, 4).Value = 44516
$ws.Cells.Item(448, 4).NumberFormat = $ws.Cells.Item(449, 4).NumberFormat
$ws.Cells.Item(448, 5).Value = 4
$ws.Cells.Item(448, 6).Value = 100112004
$ws.Cells.Item(448, 7).Value = "Cebolla"
$ws.Cells.Item(448, 8).Value = "Sin especificar"
$ws.Cells.Item(448, 9).Value = "2a (cosecha)"
$ws.Cells.Item(448, 10).Value = 1680
$ws.Cells.Item(448, 11).Value = 4500
$ws.Cells.Item(448, 12).Value = 4600
$ws.Cells.Item(448, 13).Value = 4550
$ws.Cells.Item(448, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(448, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(448, 16).Value = 253
$ws.Cells.Item(448, 17).Value = 18
$ws.Cells.Item(448, 18).Value = "Hortaliza"
